# Update crypto price/volume figures per the Dec 2 2023 GitHub Actions refresh.
# Values are written with a leading apostrophe so Excel keeps them as literal text
# (matching the inline-string cells already in the sheet) instead of auto-parsing
# numeric-looking strings (e.g. "228.63") into real numbers. The Style reset that
# follows clears the "quote prefix" formatting flag Excel attaches when it does this,
# so no new cell style is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = '''38.747.07'
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.Value = '''  +1.00%  '
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.Value = '''2.098.17'
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.Value = '''  -0.17%  '
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.Value = '''  +0.03%  '
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.Value = '''228.63'
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.Value = '''  -0.21%  '
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.Value = '''  +0.59%  '
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.Value = '''62.24'
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.Value = '''  +1.38%  '
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.Value = '''  -0.02%  '
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.Value = '''  +1.65%  '
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.Value = '''  -0.68%  '
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.Value = '''  +0.31%  '
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.Value = '''15.84'
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.Value = '''  +7.24%  '
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.Value = '''2.410.74'
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.Value = '''  -0.12%  '
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.Value = '''22.15'
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.Value = '''  -1.04%  '
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.Value = '''  +3.43%  '
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.Value = '''5.50'
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.Value = '''  +0.26%  '
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.Value = '''2.094.19'
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.Value = '''  -1.24%  '
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.Value = '''38.739.53'
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.Value = '''  +1.23%  '
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.Value = '''72.01'
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.Value = '''  +2.25%  '
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.Value = '''6.04'
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.Value = '''  +0.26%  '
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.Value = '''0.0₃0838'
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.Value = '''  +0.41%  '
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.Value = '''227.87'
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.Value = '''  +1.51%  '
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.Value = '''  +0.00%  '
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.Value = '''2.36'
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.Value = '''  -3.18%  '
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.Value = '''  +0.40%  '
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.Value = '''171.97'
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.Value = '''  +1.25%  '
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.Value = '''9.58'
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.Value = '''  +1.58%  '
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.Value = '''0.140'
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.Value = '''  +6.86%  '
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.Value = '''  +4.09%  '
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.Value = '''19.33'
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.Value = '''  +1.53%  '
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.Value = '''  +3.47%  '
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.Value = '''  +0.73%  '
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.Value = '''  +1.95%  '
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.Value = '''4.76'
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.Value = '''  +0.73%  '
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.Value = '''  +2.32%  '
$cell.Style = "Normal"
$cell = $ws.Range("D36")
$cell.Value = '''6.58'
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.Value = '''  +2.67%  '
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.Value = '''  +0.66%  '
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.Value = '''  +0.16%  '
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.Value = '''18.31'
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.Value = '''  +0.98%  '
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.Value = '''  +4.01%  '
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.Value = '''102.20'
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.Value = '''  +2.17%  '
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.Value = '''1.532.62'
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.Value = '''  -1.10%  '
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.Value = '''  -0.90%  '
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.Value = '''  +3.37%  '
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.Value = '''  +0.16%  '
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.Value = '''  +1.87%  '
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.Value = '''4.14'
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.Value = '''  -0.62%  '
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.Value = '''  +1.10%  '
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.Value = '''  -1.10%  '
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.Value = '''2.297.21'
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.Value = '''  -0.10%  '
$cell.Style = "Normal"
